# Rename the inline-picture identifiers for the three logo drawings in
# this document's headers/footers:
#   - footer (wp:docPr id="3") Pearson logo : image1.png -> image2.png
#   - footer (wp:docPr id="2") Pearson logo : image1.png -> image2.png
#   - header (wp:docPr id="1") BTEC logo    : image2.jpg -> image1.jpg
# Each <wp:docPr .../> is immediately followed (inside the same
# <w:drawing>) by a matching <pic:cNvPr .../> that carries the same
# "name" attribute and must be renamed in lock-step.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

function Replace-AttrInTagAt($text, $anchor, $oldAttr, $newAttr) {
    $idx = $text.IndexOf($anchor)
    if ($idx -lt 0) { throw "Anchor not found: $anchor" }
    $tagEnd = $text.IndexOf("/>", $idx)
    if ($tagEnd -lt 0) { throw "Tag end not found after anchor" }
    $tagEnd = $tagEnd + 2
    $tag = $text.Substring($idx, $tagEnd - $idx)
    if ($tag.IndexOf($oldAttr) -lt 0) { throw "Attribute '$oldAttr' not found in tag: $tag" }
    $newTag = $tag.Replace($oldAttr, $newAttr)
    return $text.Substring(0, $idx) + $newTag + $text.Substring($tagEnd)
}

function Rename-LogoPair($xml, $docPrAnchor, $cNvPrAnchorTemplate, $oldName, $newName) {
    $startIdx = $xml.IndexOf($docPrAnchor)
    if ($startIdx -lt 0) { throw "docPr anchor not found: $docPrAnchor" }

    $xml = Replace-AttrInTagAt $xml $docPrAnchor ('name="' + $oldName + '"') ('name="' + $newName + '"')

    $cNvPrAnchor = $cNvPrAnchorTemplate
    $cNvIdx = $xml.IndexOf($cNvPrAnchor, $startIdx)
    if ($cNvIdx -lt 0) { throw "pic:cNvPr anchor not found after docPr: $cNvPrAnchor" }
    if ($cNvIdx - $startIdx -gt 2000) { throw "pic:cNvPr match too far from docPr anchor" }

    $xml = Replace-AttrInTagAt $xml $cNvPrAnchor ('name="' + $oldName + '"') ('name="' + $newName + '"')

    return $xml
}

# --- Footer (first): Pearson logo, wp:docPr id="3" ---
$xml = Rename-LogoPair $xml `
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/>' `
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>' `
    'image1.png' 'image2.png'

# --- Footer (second): Pearson logo, wp:docPr id="2" ---
$xml = Rename-LogoPair $xml `
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/>' `
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>' `
    'image1.png' 'image2.png'

# --- Header: BTEC logo, wp:docPr id="1" ---
$xml = Rename-LogoPair $xml `
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>' `
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>' `
    'image2.jpg' 'image1.jpg'

$d.WordOpenXML = $xml

Write-Output "Renamed 3 logo picture pairs"
